# Scheduled-runner refresh of the "Famfrit_Profits" price/profit sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Updates currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H:N) for the leves whose underlying
# market data moved since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 146.26666
$ws.Range("I6").Value = 167.46153
$ws.Range("J6").Value = 8.5
$ws.Range("K6").Value = 502.38459
$ws.Range("L6").Value = 25.5
$ws.Range("M6").Value = -390.38459
$ws.Range("N6").Value = -249.5

$ws.Range("H15").Value = 674.01697
$ws.Range("I15").Value = 674.01697
$ws.Range("J15").Value = 0.0
$ws.Range("K15").Value = 2022.05091
$ws.Range("L15").Value = 0.0
$ws.Range("M15").Value = -1853.05091

$ws.Range("H17").Value = 11116146.0
$ws.Range("I17").Value = 0.0
$ws.Range("J17").Value = 11116146.0
$ws.Range("K17").Value = 0.0
$ws.Range("L17").Value = 33348438.0
$ws.Range("N17").Value = -33348774.0

$ws.Range("H28").Value = 1693.579
$ws.Range("I28").Value = 489.1111
$ws.Range("J28").Value = 2777.6
$ws.Range("K28").Value = 489.1111
$ws.Range("L28").Value = 2777.6
$ws.Range("M28").Value = -4.111100000000022
$ws.Range("N28").Value = -3747.6

$ws.Range("H40").Value = 3335383.2
$ws.Range("I40").Value = 5001000.0
$ws.Range("J40").Value = 4150.0
$ws.Range("K40").Value = 5001000.0
$ws.Range("L40").Value = 4150.0
$ws.Range("M40").Value = -5000825.0
$ws.Range("N40").Value = -4500.0

$ws.Range("H53").Value = 1816.9
$ws.Range("I53").Value = 1562.3334
$ws.Range("J53").Value = 2198.75
$ws.Range("K53").Value = 1562.3334
$ws.Range("L53").Value = 2198.75
$ws.Range("M53").Value = -925.3334
$ws.Range("N53").Value = -3472.75

$ws.Range("H113").Value = 3071.2856
$ws.Range("I113").Value = 2099.6667
$ws.Range("K113").Value = 2099.6667
$ws.Range("M113").Value = 1154.3333

$ws.Range("H116").Value = 15749.25
$ws.Range("I116").Value = 15749.25
$ws.Range("K116").Value = 15749.25
$ws.Range("M116").Value = -12307.25

$ws.Range("H132").Value = 2026.6666
$ws.Range("I132").Value = 1934.6666
$ws.Range("K132").Value = 5803.9998
$ws.Range("M132").Value = -3273.9998

$ws.Range("H135").Value = 727.25714
$ws.Range("I135").Value = 490.48
$ws.Range("K135").Value = 4414.32
$ws.Range("M135").Value = -1879.32

$ws.Range("H138").Value = 6806736.5
$ws.Range("I138").Value = 780.5263
$ws.Range("J138").Value = 11117175.0
$ws.Range("K138").Value = 2341.5789
$ws.Range("L138").Value = 33351525.0
$ws.Range("M138").Value = 2798.4211
$ws.Range("N138").Value = -33361805.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 71428760.0
$ws.Range("I5").Value = 200.84616
$ws.Range("J5").Value = 1000000000.0
$ws.Range("K5").Value = 200.84616
$ws.Range("L5").Value = 1000000000.0
$ws.Range("M5").Value = -88.84616
$ws.Range("N5").Value = -1000000224.0

$ws.Range("H32").Value = 6083.838
$ws.Range("I32").Value = 4784.5225
$ws.Range("K32").Value = 4784.5225
$ws.Range("M32").Value = -4497.5225

$ws.Range("H34").Value = 139600.0
$ws.Range("J34").Value = 265000.0
$ws.Range("L34").Value = 265000.0
$ws.Range("N34").Value = -265542.0

$ws.Range("H132").Value = 125202300.0
$ws.Range("I132").Value = 36409.332
$ws.Range("K132").Value = 109227.996
$ws.Range("M132").Value = -106697.996

$ws.Range("H138").Value = 60000.0
$ws.Range("J138").Value = 60000.0
$ws.Range("L138").Value = 60000.0
$ws.Range("N138").Value = -70280.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 71428760.0
$ws.Range("I4").Value = 200.84616
$ws.Range("J4").Value = 1000000000.0
$ws.Range("K4").Value = 200.84616
$ws.Range("L4").Value = 1000000000.0
$ws.Range("M4").Value = -85.84616
$ws.Range("N4").Value = -1000000230.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3728.825
$ws.Range("I31").Value = 2803.9524
$ws.Range("K31").Value = 2803.9524
$ws.Range("M31").Value = -2508.9524

$ws.Range("H34").Value = 3728.825
$ws.Range("I34").Value = 2803.9524
$ws.Range("K34").Value = 2803.9524
$ws.Range("M34").Value = -2601.9524

$ws.Range("H35").Value = 9999.0
$ws.Range("J35").Value = 0.0
$ws.Range("L35").Value = 0.0
$ws.Range("N35").ClearContents()

$ws.Range("H69").Value = 76900.0
$ws.Range("I69").Value = 88375.0
$ws.Range("J69").Value = 31000.0
$ws.Range("K69").Value = 88375.0
$ws.Range("L69").Value = 31000.0
$ws.Range("M69").Value = -87626.0
$ws.Range("N69").Value = -32498.0

$ws.Range("H72").Value = 76900.0
$ws.Range("I72").Value = 88375.0
$ws.Range("J72").Value = 31000.0
$ws.Range("K72").Value = 265125.0
$ws.Range("L72").Value = 93000.0
$ws.Range("M72").Value = -261381.0
$ws.Range("N72").Value = -100488.0

$ws.Range("H132").Value = 45816.824
$ws.Range("I132").Value = 60589.03
$ws.Range("J132").Value = 3962.25
$ws.Range("K132").Value = 181767.09
$ws.Range("L132").Value = 11886.75
$ws.Range("M132").Value = -179237.09
$ws.Range("N132").Value = -16946.75

$ws.Range("H134").Value = 1216.1111
$ws.Range("I134").Value = 706.4286
$ws.Range("J134").Value = 3000.0
$ws.Range("K134").Value = 2119.2858
$ws.Range("L134").Value = 9000.0
$ws.Range("M134").Value = 415.7142000000003
$ws.Range("N134").Value = -14070.0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27210752.0
$ws.Range("I4").Value = 34497330.0
$ws.Range("J4").Value = 17148334.0
$ws.Range("K4").Value = 103491990.0
$ws.Range("L4").Value = 51445002.0
$ws.Range("M4").Value = -103491878.0
$ws.Range("N4").Value = -51445226.0

$ws.Range("H55").Value = 15152188.0
$ws.Range("J55").Value = 33334832.0
$ws.Range("L55").Value = 100004496.0
$ws.Range("N55").Value = -100004850.0

$ws.Range("H130").Value = 2699.2
$ws.Range("I130").Value = 1030.0
$ws.Range("J130").Value = 3116.5
$ws.Range("K130").Value = 3090.0
$ws.Range("L130").Value = 9349.5
$ws.Range("M130").Value = 1930.0
$ws.Range("N130").Value = -19389.5

$ws.Range("H140").Value = 1245.5862
$ws.Range("I140").Value = 779.4783
$ws.Range("J140").Value = 3032.3333
$ws.Range("K140").Value = 2338.4349
$ws.Range("L140").Value = 9096.999899999999
$ws.Range("M140").Value = 2841.5651
$ws.Range("N140").Value = -19456.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 254333.88
$ws.Range("I70").Value = 403940.8
$ws.Range("J70").Value = 4989.0
$ws.Range("K70").Value = 403940.8
$ws.Range("L70").Value = 4989.0
$ws.Range("M70").Value = -403670.8
$ws.Range("N70").Value = -5529.0

$ws.Range("H73").Value = 254333.88
$ws.Range("I73").Value = 403940.8
$ws.Range("J73").Value = 4989.0
$ws.Range("K73").Value = 403940.8
$ws.Range("L73").Value = 4989.0
$ws.Range("M73").Value = -403004.8
$ws.Range("N73").Value = -6861.0

$ws.Range("H97").Value = 2350.2964
$ws.Range("I97").Value = 1675.8948
$ws.Range("J97").Value = 3952.0
$ws.Range("K97").Value = 1675.8948
$ws.Range("L97").Value = 3952.0
$ws.Range("M97").Value = -1179.8948
$ws.Range("N97").Value = -4944.0

$ws.Range("H102").Value = 2394.3794
$ws.Range("I102").Value = 1330.1111
$ws.Range("K102").Value = 1330.1111
$ws.Range("M102").Value = 291.8888999999999

$ws.Range("H107").Value = 664.53845
$ws.Range("I107").Value = 561.3333
$ws.Range("K107").Value = 561.3333
$ws.Range("M107").Value = 1358.6667

$ws.Range("H132").Value = 2585.7144
$ws.Range("I132").Value = 2585.7144
$ws.Range("J132").Value = 0.0
$ws.Range("K132").Value = 7757.1432
$ws.Range("L132").Value = 0.0
$ws.Range("M132").Value = -5227.1432
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3891.5
$ws.Range("I40").Value = 3669.8
$ws.Range("K40").Value = 3669.8
$ws.Range("M40").Value = -3533.8

$ws.Range("H46").Value = 1186.56
$ws.Range("I46").Value = 669.2857
$ws.Range("K46").Value = 669.2857
$ws.Range("M46").Value = -481.2857

$ws.Range("H55").Value = 617.2083
$ws.Range("I55").Value = 423.57144
$ws.Range("J55").Value = 888.3
$ws.Range("K55").Value = 423.57144
$ws.Range("L55").Value = 888.3
$ws.Range("M55").Value = -250.57144
$ws.Range("N55").Value = -1234.3

$ws.Range("H122").Value = 4613.2646
$ws.Range("I122").Value = 3654.7144
$ws.Range("K122").Value = 10964.1432
$ws.Range("M122").Value = -8514.143199999999

$ws.Range("H132").Value = 11414.857
$ws.Range("I132").Value = 9555.192
$ws.Range("K132").Value = 28665.576
$ws.Range("M132").Value = -26135.576

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 69000.0
$ws.Range("J114").Value = 69000.0
$ws.Range("L114").Value = 69000.0
$ws.Range("N114").Value = -77678.0

$ws.Range("H123").Value = 59999.332
$ws.Range("J123").Value = 59999.332
$ws.Range("L123").Value = 59999.332
$ws.Range("N123").Value = -69799.332

$ws.Range("H132").Value = 2328.641
$ws.Range("I132").Value = 2786.724
$ws.Range("K132").Value = 8360.172
$ws.Range("M132").Value = -5830.172

$ws.Range("H139").Value = 49000.0
$ws.Range("J139").Value = 49000.0
$ws.Range("L139").Value = 49000.0
$ws.Range("N139").Value = -59280.0
